$wb = $excel.ActiveWorkbook

function Set-Row($ws, $row, $a, $b, $c, $d, $e) {
  $ws.Cells.Item($row,1).Value = $a
  $ws.Cells.Item($row,2).Value = $b
  $ws.Cells.Item($row,3).Value = $c
  $ws.Cells.Item($row,4).Value = $d
  $ws.Cells.Item($row,5).Value = $e
}

# ----------------------------------------------------------------------
# Sheet 1: LP1912  (70 -> 79 data rows; new scrape run at 07:48:31)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

Set-Row $ws1 76 "07:48:31" "07:48" "14_ABASTO"          0   "LP1912"
Set-Row $ws1 77 "07:48:31" "08:10" "16_SANTA ANA"       22  "LP1912"
Set-Row $ws1 78 "07:48:31" "08:32" "23_HERNANDEZ"       44  "LP1912"
Set-Row $ws1 79 "07:48:31" "08:53" "10_OLMOS"           65  "LP1912"
Set-Row $ws1 80 "07:48:31" "09:07" "23_HERNANDEZ"       79  "LP1912"
Set-Row $ws1 81 "07:48:31" "09:23" "11_ETCHEVERRY"      95  "LP1912"
Set-Row $ws1 82 "07:48:31" "09:32" "15_ABASTO"          104 "LP1912"
Set-Row $ws1 83 "07:48:31" "09:33" "10_OLMOS"           105 "LP1912"
Set-Row $ws1 84 "07:48:31" "09:42" "215C_EL PATO"       114 "LP1912"

$rng1 = $ws1.Range("A6:E84")
$key1 = $ws1.Range("B6:B84")
$rng1.Sort($key1)

$ws1.Range("A2").Value = "Última actualización: 07:48:31"
$ws1.Range("A3").Value = "Total filas: 79"

# ----------------------------------------------------------------------
# Sheet 2: LP1912-215  (14 -> 15 data rows)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

Set-Row $ws2 20 "07:48:31" "09:42" "215C_EL PATO" 114 "LP1912"

$rng2 = $ws2.Range("A6:E20")
$key2 = $ws2.Range("B6:B20")
$rng2.Sort($key2)

$ws2.Range("A2").Value = "Última actualización: 07:48:31"
$ws2.Range("A3").Value = "Total filas: 15"

# ----------------------------------------------------------------------
# Sheet 3: 6203-6173  (17 -> 19 data rows)
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

Set-Row $ws3 23 "07:48:31" "08:36" "215A_LA PLATA" 48 "L6173"
Set-Row $ws3 24 "07:48:31" "09:09" "215D_LA PLATA" 81 "L6203"

$rng3 = $ws3.Range("A6:E24")
$key3 = $ws3.Range("B6:B24")
$rng3.Sort($key3)

$ws3.Range("A2").Value = "Última actualización: 07:48:31"
$ws3.Range("A3").Value = "Total filas: 19"
